# "couple quick edits to each table"
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Loop Table 1")
$ws2 = $wb.Worksheets.Item("Loop Table 2")

# ---------------------------------------------------------------------------
# Loop Table 1: the old loop stopped at x=16 (<=17); update the story so the
# loop now goes one more round and stops at x=18 (>17), adding a new row.
# ---------------------------------------------------------------------------

# Drop the old explanatory note in E15 first so the shared-string table can
# garbage-collect it before the new strings are appended (keeps the shared
# string order/index identical to the target).
$ws1.Range("E15").ClearContents()

# New row 16: loop count 10, value of x = 18, with the new explanation.
$b16 = $ws1.Range("B16")
$b16.Value = 10
$b16.HorizontalAlignment = -4108
$b16.VerticalAlignment = -4108

$c16 = $ws1.Range("C16")
$c16.Value = 18
$c16.HorizontalAlignment = -4108
$c16.Font.Bold = $true

$ws1.Range("E16").Value = "the loop stops here, because 18 is not less than or equal to 17"

# ---------------------------------------------------------------------------
# Loop Table 2: add the final loop iteration (index 8, which ends the loop)
# plus a result row showing the final concatenation, merged across E:F.
# ---------------------------------------------------------------------------

$b22 = $ws2.Range("B22")
$b22.Value = 9
$b22.HorizontalAlignment = -4108
$b22.VerticalAlignment = -4108

$c22 = $ws2.Range("C22")
$c22.Value = 8
$c22.HorizontalAlignment = -4108
$c22.VerticalAlignment = -4108
$c22.Font.Bold = $true

$b24c24 = $ws2.Range("B24:C24")
$b24c24.HorizontalAlignment = -4108
$b24c24.HorizontalAlignment = 1

$e24f24 = $ws2.Range("E24:F24")
$e24f24.HorizontalAlignment = -4108
$e24f24.Value = "meowpurr"
$e24f24.Merge()

# ---------------------------------------------------------------------------
# Selection / active-sheet bookkeeping to match the saved view state.
# ---------------------------------------------------------------------------

$ws2.Select()
$ws2.Range("G24").Select()

$ws1.Select()
$ws1.Range("C16").Select()
